# Powerpoint writer: consolidate text run nodes.
# Merge every "word" run that is immediately followed by a lone
# single-space run into one run ("word" + " "), trimming the
# total run count without changing the visible text.

function Merge-Runs {
    param($TextRange, [string[]]$Parts)

    $pos = 1
    foreach ($part in $Parts) {
        $len = $part.Length
        $TextRange.Characters($pos, $len).Text = $part
        $pos = $pos + $len
    }
}

$p = $ppt.ActivePresentation

$titlePlans = @{
    1  = @("Slide ", "1 ", "(Content)")
    2  = @("Slide ", "2 ", "(Content)")
    3  = @("Slide ", "3 ", "(Content)")
    4  = @("Slide ", "4 ", "(Content)")
    5  = @("Slide ", "5 ", "(Two ", "Content)")
    6  = @("Slide ", "6 ", "(Two ", "Content ", "Right)")
    7  = @("Slide ", "7 ", "(Content ", "with ", "Caption)")
    8  = @("Slide ", "8 ", "(Comparison)")
    9  = @("Slide ", "9 ", "(Content)")
    10 = @("Slide ", "10 ", "(Content)")
    11 = @("Slide ", "11 ", "(Content)")
    12 = @("Slide ", "12 ", "(Content)")
}

foreach ($idx in 1..12) {
    $slide = $p.Slides.Item($idx)
    $title = $slide.Shapes.Item(1)
    Merge-Runs $title.TextFrame.TextRange $titlePlans[$idx]
}

# "an image" / "An image" captions sit in "TextBox 3" on slides 6-8.
$captionPlan = @{
    6 = "an "
    7 = "An "
    8 = "An "
}

foreach ($idx in 6..8) {
    $slide = $p.Slides.Item($idx)
    $box = $slide.Shapes.Item("TextBox 3")
    Merge-Runs $box.TextFrame.TextRange @($captionPlan[$idx], "image")
}
